$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header/data cells - column B filled first, then column C
# (matches the shared-string insertion order seen in the target workbook)
$ws.Range("B1").Value = "Producto"
$ws.Range("B2").Value = "Cerveza Club Colombia Dorada lata x6und x330ml c-u"
$ws.Range("C1").Value = "CorreoAsociarCompra"
$ws.Range("C2").Value = "paangudi3@gmail.com"

# Turn the e-mail address into a mailto hyperlink (this also creates the
# builtin "Hyperlink" cell style used by C2)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:paangudi3@gmail.com")

# Size the new columns to fit their content
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()

# Move the active selection like in the final workbook
$ws.Range("C3").Select()
